$p = $ppt.ActivePresentation

# --- Slide 13 gets a new "Text Placeholder 2" (body idx=10) shape.
#     Slide 14 already owns exactly that placeholder, so copy it over. ---
$srcShape = $p.Slides.Item(14).Shapes.Item(2)
$srcShape.Copy()
$p.Slides.Item(13).Shapes.Paste() | Out-Null

# --- Give every content slide (2-25) its title text "Slide N". ---
for ($i = 2; $i -le 25; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = "Slide " + $i
}
